# Add "generic_text_elements" and "case_text_elements" input sheets,
# mirroring the header layout/style already used by the other
# "label" + "value" sheets in this template (e.g. key_output_weights).

$wb = $excel.ActiveWorkbook

$cfgSheet = $wb.Worksheets.Item("configurations")
$headerSrc = $wb.Worksheets.Item("key_output_weights")

# --- create both new sheets in their final tab order first ---
# generic_text_elements right after "configurations" ...
$genericSheet = $wb.Worksheets.Add($null, $cfgSheet)
$genericSheet.Name = "generic_text_elements"

# ... and case_text_elements right after that.
$caseSheet = $wb.Worksheets.Add($null, $genericSheet)
$caseSheet.Name = "case_text_elements"

# Fill in case_text_elements' header first, then generic_text_elements',
# so new shared-string entries are interned in the same order as the
# source edit (case_text_element, then generic_text_element).
$caseSheet.Range("A1").Value = "case_text_element"
$caseSheet.Range("B1").Value = "value"

$genericSheet.Range("A1").Value = "generic_text_element"
$genericSheet.Range("B1").Value = "value"

# Reuse the same header formatting (bold font, fill, border, center/top
# alignment) already defined for the other sheets' header rows.
$headerSrc.Range("A1:B1").Copy() | Out-Null
$genericSheet.Range("A1:B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$headerSrc.Range("A1:B1").Copy() | Out-Null
$caseSheet.Range("A1:B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$genericSheet.Columns.Item(1).ColumnWidth = 17.833333333333332
$genericSheet.Columns.Item(2).ColumnWidth = 4.5

$caseSheet.Columns.Item(1).ColumnWidth = 15.166666666666666
$caseSheet.Columns.Item(2).ColumnWidth = 4.5

$excel.CutCopyMode = 0

# Make "generic_text_elements" the active tab/selection, matching the
# new workbook view, and set a realistic last-used selection per sheet.
$genericSheet.Activate()
$genericSheet.Range("G33").Select() | Out-Null

$caseSheet.Range("C6").Select() | Out-Null

$wb.Worksheets.Item("key_outputs").Range("F5").Select() | Out-Null

$genericSheet.Activate()
